$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B11: value changes to the text "1".
# It must remain a *text* cell (shared string), not get auto-converted to a
# number, and it must keep its original style (border/fill) untouched.
# A plain Value="1" assignment gets coerced to a number, and forcing text via
# a leading apostrophe / NumberFormat="@" creates a brand-new style entry
# (quotePrefix / custom number format) instead of reusing the existing one.
# Routing the literal through a TEXT() formula and then hardening it back to
# a plain value in-place (copy + paste-special values-only) keeps the
# original cell style and produces a genuine text value.
$ws.Range("B11").Formula = "=TEXT(1,""0"")"
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
